$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.877.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0621"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.901.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.665.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.550"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "249.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.833.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0732"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +6.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.432.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.43%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.49%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.930"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.583"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("E44").Value = "  -6.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.809.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.787"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0111"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("E51").Value = "  -5.01%  "
